$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers I1 (I0) and J1 (IF), copying H1's header formatting
# (bold font, thin border, centered/top alignment) so the new header cells
# match the style of the existing headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for I2:I87 and J2:J87 (row 2 corresponds to index 0, etc.)
$iVals = @(8,8,9,9,7,9,9,7,9,8,9,8,9,8,10,9,9,8,9,9,9,9,10,9,8,8,9,9,9,10,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,10,9,8,9,10,9,9,7,9,8,9,8,7,8,10,9,8,9,9,8,4,7,7,7,5,4,4,3)
$jVals = @(8,8,9,9,8,9,9,7,9,8,9,8,9,8,10,9,9,8,9,9,9,9,10,9,8,8,9,9,9,10,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,10,9,8,9,10,9,9,7,9,8,9,8,7,8,10,9,8,9,9,8,4,7,7,7,5,4,4,3)

for ($r = 2; $r -le 87; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
